$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 211,
# pushing all subsequent rows (211-247) down by one (212-248).
$ws.Rows(211).Insert()

# Populate the newly inserted row 211 with its data.
$ws.Range("A211").Value = 3
$ws.Range("B211").Value = "Femacal de La Calera"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44522
$ws.Range("E211").Value = 5
$ws.Range("F211").Value = 100112040
$ws.Range("G211").Value = "Cilantro"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 60
$ws.Range("K211").Value = 3000
$ws.Range("L211").Value = 3000
$ws.Range("M211").Value = 3000
$ws.Range("N211").Value = "$/docena de atados (3 kilos)"
$ws.Range("O211").Value = "Provincia de Quillota"
$ws.Range("P211").Value = 1000
$ws.Range("Q211").Value = 3
$ws.Range("R211").Value = "Hortaliza"
